$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("R2").Value = 1.41
$ws.Range("R3").Value = 1.41
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("R4").Value = 1.67
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("Q5").Value = 2.6
$ws.Range("R5").Value = 1.44
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10
$ws.Range("Q7").Value = 2.08
$ws.Range("R7").Value = 1.73
$ws.Range("Q8").Value = 2.2
$ws.Range("R8").Value = 1.65
$ws.Range("G14").Value = 1.8
$ws.Range("I14").Value = 4.75
$ws.Range("J14").Value = 2.4
$ws.Range("Z14").Value = 15
$ws.Range("AK14").Value = 34
$ws.Range("AO14").Value = 9.5
$ws.Range("AW14").Value = 6
$ws.Range("O19").Value = 1.67
$ws.Range("P19").Value = 2.1
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = 7.9
$ws.Range("I20").Value = 1.07
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = 3.35
$ws.Range("L20").Value = 1.33
$ws.Range("O20").Value = 1.08
$ws.Range("P20").Value = 6.5
$ws.Range("Q20").Value = 1.27
$ws.Range("R20").Value = 3.45
$ws.Range("S20").Value = 1.18
$ws.Range("T20").Value = 4.25
$ws.Range("U20").Value = 2.35
$ws.Range("V20").Value = 1.53
$ws.Range("W20").Value = 120
$ws.Range("X20").Value = 500
$ws.Range("Y20").Value = 120
$ws.Range("AA20").Value = 700
$ws.Range("AB20").Value = 350
$ws.Range("AD20").Value = 23
$ws.Range("AE20").Value = 50
$ws.Range("AF20").Value = 200
$ws.Range("AG20").Value = 11.75
$ws.Range("AH20").Value = 7.1
$ws.Range("AI20").Value = 13.5
$ws.Range("AK20").Value = 11.5
$ws.Range("AL20").Value = 40
$ws.Range("AN20").Value = 23
$ws.Range("AO20").Value = 200
$ws.Range("AP20").Value = 100
$ws.Range("AT20").Value = 4.25
$ws.Range("AU20").Value = 11.75
$ws.Range("AV20").Value = 100
$ws.Range("AW20").Value = 3.15
$ws.Range("AX20").Value = 4.05
$ws.Range("AZ20").Value = 7.5
$ws.Range("BA20").Value = 28
